# Generate Report for Handoff
# Renames the handed-off markdown file's GUID (and dependent xlf hashes),
# refreshing the "Overview"/"zh-cn"/"de-de" sheets of the localization
# status report, and bumps the handoff timestamps forward.

$wb = $excel.ActiveWorkbook

$newGuid = "8bfb0a3a-5cd2-455c-a925-2183d4c922a7"
$newHash = "bc98418451f317f8dbd78842b4cfa79c9a898c78"

$newMdName      = "$newGuid.md"
$newMdPath      = "e2e\$newGuid.md"
$newZhXlfName   = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName   = "$newGuid.$newHash.de-de.xlf"

$newOverviewDate = "2016-08-27 10:57:10"
$newHandoffDate  = "2016-08-27 10:57:05"

# The hyperlink target (commit blob URL) itself is not changing - only the
# display text shown in the cell changes - so reuse the existing address.
$hlAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c8d0813b11bb9b60da0cbb2a5e379d8a5a02d961/e2e/1ce8c6e3-1dde-4dd0-b5fa-cb5c81809e61.md"

# ---------------------------------------------------------------------
# Overview sheet: File Name (A2), Path And Name (B2, hyperlink),
# Latest HO Xliff Generate Date (G2)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hlAddress, [System.Type]::Missing, [System.Type]::Missing, $newMdPath) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: Source File Name (A2, hyperlink), Latest Handoff File
# (G2), Latest Handoff Datetime (H2)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("G2").Value = $newZhXlfName
$wsZh.Range("H2").Value = $newHandoffDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hlAddress, [System.Type]::Missing, [System.Type]::Missing, $newMdName) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: Source File Name (A2, hyperlink), Latest Handoff File
# (G2), Latest Handoff Datetime (H2 - shares the Overview timestamp)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("G2").Value = $newDeXlfName
$wsDe.Range("H2").Value = $newOverviewDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hlAddress, [System.Type]::Missing, [System.Type]::Missing, $newMdName) | Out-Null

Write-Host "Handoff report regenerated for $newGuid"
